$p = $ppt.ActivePresentation

# Remove slides 6 through 21 (sldId 261-276 / rId7-rId22),
# keeping only the first 5 slides.
for ($i = $p.Slides.Count; $i -ge 6; $i--) {
    $p.Slides.Item($i).Delete()
}
